$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SamplesTab" row (row 3) query previously included the
# sample_tumor_status / sample_type columns ("Tumor" / "Analyte Type").
# The new "CDS All studies" testcase query drops those two extra
# columns, keeping just Sample ID / Participant ID / Study Name /
# Accession.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND smp.sample_type = 'Blood Derived Normal'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesQuery

# Move the active selection up one row (was C4, now C3), matching the
# refreshed view state saved with the workbook.
$ws.Range("C3").Select() | Out-Null
